# "17 Apr 25 Changes, work on well schematic window."
# Update the checklist on the "3. Well Builder" sheet:
#  - rename the row 14 task text
#  - mark rows 10-13 (and new rows 15/16) as done ("X" in column A)
#  - add new to-do rows 15-21 describing further well-schematic work

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3. Well Builder")
$ws.Activate()

# Row 9 keeps its text ("Change casing") - nothing to change there.

# Mark existing rows 10-13 as completed.
$ws.Range("A10").Value2 = "X"
$ws.Range("A11").Value2 = "X"
$ws.Range("A12").Value2 = "X"
$ws.Range("A13").Value2 = "X"

# Row 14: rename the outstanding task.
$ws.Range("C14").Value2 = "Update casing info window based on casing in database"

# New row 15: completed task about updating the drawing.
$ws.Range("A15").Value2 = "X"
$ws.Range("B15").Value2 = 7
$ws.Range("C15").Value2 = "update drawing after adding casing"

# New row 16: completed sub-task under 7 (7.1).
$ws.Range("A16").Value2 = "X"
$ws.Range("C16").Value2 = 7.1
$ws.Range("D16").Value2 = "clear canvas somehow?"

# New rows 17-21: outstanding to-do items.
$ws.Range("B17").Value2 = 8
$ws.Range("C17").Value2 = "why does it say a casing already exists"

$ws.Range("B18").Value2 = 9
$ws.Range("C18").Value2 = "Add conductor"

$ws.Range("B19").Value2 = 10
$ws.Range("C19").Value2 = "find way to order casing correctly in database (surf->int->liner->prod->tubing)"

$ws.Range("B20").Value2 = 11
$ws.Range("C20").Value2 = "draw packer with liners and tubing"

$ws.Range("B21").Value2 = 12
$ws.Range("C21").Value2 = "Show formation depths"

# Update the selection to match the state after these edits.
[void]$ws.Range("B22").Select()
